$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1558.8889
$ws.Range("I6").Value = 1558.8889
$ws.Range("K6").Value = 4676.6667
$ws.Range("M6").Value = -4564.6667
$ws.Range("H17").Value = 6699999.5
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H33").Value = 368.375
$ws.Range("I33").Value = 227
$ws.Range("J33").Value = 729.6667
$ws.Range("K33").Value = 227
$ws.Range("L33").Value = 729.6667
$ws.Range("M33").Value = 2
$ws.Range("N33").Value = -1187.6667
$ws.Range("H111").Value = 2000
$ws.Range("I111").Value = 2000
$ws.Range("K111").Value = 6000
$ws.Range("M111").Value = -2933
$ws.Range("H132").Value = 3169.0566
$ws.Range("I132").Value = 3162.1956
$ws.Range("K132").Value = 9486.586800000001
$ws.Range("M132").Value = -6956.586800000001
$ws.Range("H135").Value = 1870.5
$ws.Range("I135").Value = 1359.1
$ws.Range("J135").Value = 3149
$ws.Range("K135").Value = 12231.9
$ws.Range("L135").Value = 28341
$ws.Range("M135").Value = -9696.9
$ws.Range("N135").Value = -33411
$ws.Range("H137").Value = 4453.852
$ws.Range("I137").Value = 3651.611
$ws.Range("K137").Value = 10954.833
$ws.Range("M137").Value = -8404.832999999999
$ws.Range("H138").Value = 3059.75
$ws.Range("I138").Value = 2237.2632
$ws.Range("J138").Value = 3598.6206
$ws.Range("K138").Value = 6711.7896
$ws.Range("L138").Value = 10795.8618
$ws.Range("M138").Value = -1571.7896
$ws.Range("N138").Value = -21075.8618

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12629529
$ws.Range("I32").Value = 7577021
$ws.Range("K32").Value = 7577021
$ws.Range("M32").Value = -7576734
$ws.Range("H61").Value = 2148.2424
$ws.Range("I61").Value = 1996.5161
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 1996.5161
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -1784.5161
$ws.Range("N61").Value = -4924
$ws.Range("H63").Value = 3559.8572
$ws.Range("I63").Value = 2223.5
$ws.Range("J63").Value = 5341.6665
$ws.Range("K63").Value = 2223.5
$ws.Range("L63").Value = 5341.6665
$ws.Range("M63").Value = -1537.5
$ws.Range("N63").Value = -6713.6665
$ws.Range("H66").Value = 3559.8572
$ws.Range("I66").Value = 2223.5
$ws.Range("J66").Value = 5341.6665
$ws.Range("K66").Value = 11117.5
$ws.Range("L66").Value = 26708.3325
$ws.Range("M66").Value = -7685.5
$ws.Range("N66").Value = -33572.3325
$ws.Range("H102").Value = 2568.7778
$ws.Range("I102").Value = 2394.8333
$ws.Range("K102").Value = 2394.8333
$ws.Range("M102").Value = -772.8332999999998
$ws.Range("H112").Value = 43815.734
$ws.Range("J112").Value = 43815.734
$ws.Range("L112").Value = 43815.734
$ws.Range("N112").Value = -46769.734
$ws.Range("H132").Value = 2756.9583
$ws.Range("I132").Value = 2369.861
$ws.Range("K132").Value = 7109.583
$ws.Range("M132").Value = -4579.583
$ws.Range("H136").Value = 2148.2424
$ws.Range("I136").Value = 1996.5161
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 5989.5483
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -3439.5483
$ws.Range("N136").Value = -18600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 24772
$ws.Range("I82").Value = 4748
$ws.Range("J82").Value = 64820
$ws.Range("K82").Value = 4748
$ws.Range("L82").Value = 64820
$ws.Range("M82").Value = -4365
$ws.Range("N82").Value = -65586
$ws.Range("H85").Value = 24772
$ws.Range("I85").Value = 4748
$ws.Range("J85").Value = 64820
$ws.Range("K85").Value = 4748
$ws.Range("L85").Value = 64820
$ws.Range("M85").Value = -3422
$ws.Range("N85").Value = -67472
$ws.Range("H105").Value = 3980.3076
$ws.Range("I105").Value = 3666.6667
$ws.Range("K105").Value = 3666.6667
$ws.Range("M105").Value = -1919.6667
$ws.Range("H134").Value = 13059223
$ws.Range("I134").Value = 3248459.5
$ws.Range("K134").Value = 9745378.5
$ws.Range("M134").Value = -9742843.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2362.1943
$ws.Range("I58").Value = 1876.5714
$ws.Range("K58").Value = 1876.5714
$ws.Range("M58").Value = -1673.5714
$ws.Range("H106").Value = 58446.777
$ws.Range("J106").Value = 53252.625
$ws.Range("L106").Value = 53252.625
$ws.Range("N106").Value = -55776.625
$ws.Range("H122").Value = 2476.5625
$ws.Range("I122").Value = 2342.5
$ws.Range("K122").Value = 7027.5
$ws.Range("M122").Value = -4577.5
$ws.Range("H132").Value = 1812.3334
$ws.Range("I132").Value = 1486.5333
$ws.Range("J132").Value = 3441.3333
$ws.Range("K132").Value = 4459.5999
$ws.Range("L132").Value = 10323.9999
$ws.Range("M132").Value = -1929.5999
$ws.Range("N132").Value = -15383.9999
$ws.Range("H136").Value = 2362.1943
$ws.Range("I136").Value = 1876.5714
$ws.Range("K136").Value = 5629.7142
$ws.Range("M136").Value = -3079.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1689.7142
$ws.Range("J132").Value = 1957.3
$ws.Range("L132").Value = 17615.7
$ws.Range("N132").Value = -22675.7
$ws.Range("H139").Value = 3894.1428
$ws.Range("I139").Value = 3668.1667
$ws.Range("K139").Value = 11004.5001
$ws.Range("M139").Value = -5864.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1663.1305
$ws.Range("I122").Value = 1774
$ws.Range("K122").Value = 5322
$ws.Range("M122").Value = -2872
$ws.Range("H132").Value = 2600.4443
$ws.Range("I132").Value = 2330.875
$ws.Range("J132").Value = 4757
$ws.Range("K132").Value = 6992.625
$ws.Range("L132").Value = 14271
$ws.Range("M132").Value = -4462.625
$ws.Range("N132").Value = -19331

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1905
$ws.Range("I16").Value = 1685.8
$ws.Range("J16").Value = 3001
$ws.Range("K16").Value = 1685.8
$ws.Range("L16").Value = 3001
$ws.Range("M16").Value = -1515.8
$ws.Range("N16").Value = -3341
$ws.Range("H55").Value = 278.1
$ws.Range("I55").Value = 209.875
$ws.Range("J55").Value = 356.07144
$ws.Range("K55").Value = 209.875
$ws.Range("L55").Value = 356.07144
$ws.Range("M55").Value = -36.875
$ws.Range("N55").Value = -702.0714399999999
$ws.Range("H101").Value = 55066.2
$ws.Range("J101").Value = 55066.2
$ws.Range("L101").Value = 55066.2
$ws.Range("N101").Value = -61556.2
$ws.Range("H104").Value = 23056.666
$ws.Range("J104").Value = 23056.666
$ws.Range("L104").Value = 23056.666
$ws.Range("N104").Value = -30044.666
$ws.Range("H106").Value = 1035564.2
$ws.Range("J106").Value = 1035564.2
$ws.Range("L106").Value = 1035564.2
$ws.Range("N106").Value = -1038088.2
$ws.Range("H110").Value = 72640
$ws.Range("J110").Value = 72640
$ws.Range("L110").Value = 72640
$ws.Range("N110").Value = -80820
$ws.Range("H122").Value = 3487.2
$ws.Range("I122").Value = 3487.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10461.6
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8011.599999999999
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4794.2
$ws.Range("I62").Value = 3947
$ws.Range("J62").Value = 5006
$ws.Range("K62").Value = 3947
$ws.Range("L62").Value = 5006
$ws.Range("M62").Value = -3323
$ws.Range("N62").Value = -6254
$ws.Range("H65").Value = 4794.2
$ws.Range("I65").Value = 3947
$ws.Range("J65").Value = 5006
$ws.Range("K65").Value = 19735
$ws.Range("L65").Value = 25030
$ws.Range("M65").Value = -16615
$ws.Range("N65").Value = -31270
$ws.Range("H101").Value = 58133.668
$ws.Range("J101").Value = 58133.668
$ws.Range("L101").Value = 58133.668
$ws.Range("N101").Value = -64623.668
$ws.Range("H122").Value = 3156.4736
$ws.Range("J122").Value = 5572.25
$ws.Range("L122").Value = 16716.75
$ws.Range("N122").Value = -21616.75
$ws.Range("H123").Value = 111955
$ws.Range("J123").Value = 111955
$ws.Range("L123").Value = 111955
$ws.Range("N123").Value = -121755
$ws.Range("H132").Value = 3293.36
$ws.Range("I132").Value = 2645.5715
$ws.Range("K132").Value = 7936.7145
$ws.Range("M132").Value = -5406.7145
